$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 788, shifting existing rows 788:829 down to 789:830
$ws.Rows.Item(788).Insert()

# Populate the newly inserted row 788 with its values.
# Column A holds a date-like text string; force text formatting so Excel
# doesn't auto-convert it into a real date serial number, then reset the
# style back to Normal so no stray style index is left on the cell.
$cellA = $ws.Cells.Item(788, 1)
$cellA.NumberFormat = "@"
$cellA.Value2 = "2026/02/13"
$cellA.Style = "Normal"

$ws.Cells.Item(788, 2).Value = "金"
$ws.Cells.Item(788, 3).Value = 16
$ws.Cells.Item(788, 4).Value = 201
